# Update the "About" sheet notes describing the discount rate, to reflect
# the switch from a US-centric (OMB/SCC) source to an EU-centric framing,
# as part of integrating Agora data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Row 16 (A16): reword to reference "the discount rate built into the health"
# instead of "the 3% rate used for the central estimate"
$ws.Range("A16").Value = "We choose to use a 3% discount rate here, for consistency with the discount rate built into the health"

# Row 17 (A17): used to continue the old row-16 sentence about SCoC; now
# holds the sentence that used to live in row 18.
$ws.Range("A17").Value = "damages values in the SCoHIbP Social Cost of Health Impacts by Pollutant variable."

# Row 18 is no longer used - clear it out (it previously held the sentence
# now moved to row 17).
$ws.Range("A18").ClearContents()

# Row 19 (A19): brand-new sentence about the EU's Social Cost of Carbon.
$ws.Range("A19").Value = "However, note that the EU's Social Cost of Carbon (SCoC) variable uses an estimate based on a 1% discount rate."

# Move the active selection to A20, just past the newly-added content.
$ws.Range("A20").Select() | Out-Null
